# Update condition column C: header "audioFalse" -> "currentPhase"
# and both data rows' audio-file values collapse into a single "train1P2" value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "currentPhase"
$ws.Range("C2").Value = "train1P2"
$ws.Range("C3").Value = "train1P2"
